$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the "2-3 Key Restaurants and/or Places You Enjoyed" column (D) for
# rows 5-7 with the newly captured notes. Order matters for how the shared
# string table indexes the new entries (D7, then D6, then D5).
$ws.Range("D7").Value = "Red Door Woodfired Grill, Apple, and Saint Joseph Medical Center. Visited the medical center because I injured my wrist by accidentally punching a wall in my housing. "
$ws.Range("D6").Value = "TGI Fridays, IHOP, and Texas Roadhouse. Roadhouse had literally the best steaks I had ever had - really close to Outback Steakhouse steaks."
$ws.Range("D5").Value = "Domino's Pizza, Footlocker, and Dental Associates. Two furniture stores to replace some chairs, couch cushions, and mattresses.  "

# The added text wraps within the existing styles, so the rows that now hold
# paragraph-length notes grow taller (matches the Excel auto-fit that runs on
# commit in the real workbook).
$ws.Rows.Item(6).RowHeight = 110.25
$ws.Rows.Item(7).RowHeight = 94.5

# Reflect the user's final view state: zoomed out a bit further, with B7
# selected (the workbook was also scrolled so column B became the left-most
# visible column, with row 4 near the top of the viewport).
$ws.Range("B7").Select()
$excel.ActiveWindow.Zoom = 70
